# Apply updated registration counts to the "Inscricoes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Map of cell address -> new value
$updates = @{
    "E2"  = 30
    "E4"  = 9
    "F4"  = 1
    "H4"  = 1
    "E15" = 99
    "E17" = 61
    "F17" = 25
    "H17" = 25
    "E23" = 2
    "E33" = 18
    "E36" = 51
    "E42" = 20
    "F42" = 5
    "H42" = 5
    "E44" = 16
    "E51" = 4
    "E61" = 16
    "F61" = 3
    "H61" = 3
    "E62" = 17
    "E63" = 13
    "E67" = 21
    "E68" = 8
    "E76" = 26
    "E79" = 14
    "E80" = 15
    "F81" = 2
    "H81" = 2
    "E89" = 14
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
